$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A31").Value = "First result info"
$ws.Range("A33").Value = "Second result info"
$ws.Range("A35").Value = "Third result info"
$ws.Range("A36").Value = "Navigation"

$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("B42").Select()
